$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 689.6667
$ws.Range("I29").Value = 83.166664
$ws.Range("J29").Value = 1902.6666
$ws.Range("K29").Value = 249.499992
$ws.Range("L29").Value = 5707.9998
$ws.Range("M29").Value = 31.50000800000001
$ws.Range("N29").Value = -6269.9998
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H43").Value = 964.85
$ws.Range("I43").Value = 1125.125
$ws.Range("J43").Value = 858.0
$ws.Range("K43").Value = 1125.125
$ws.Range("L43").Value = 858.0
$ws.Range("M43").Value = -1056.125
$ws.Range("N43").Value = -996.0
$ws.Range("H58").Value = 1830.8
$ws.Range("I58").Value = 1788.5
$ws.Range("J58").Value = 2000.0
$ws.Range("K58").Value = 5365.5
$ws.Range("L58").Value = 6000.0
$ws.Range("M58").Value = -5215.5
$ws.Range("N58").Value = -6300.0
$ws.Range("H61").Value = 262.5
$ws.Range("I61").Value = 262.5
$ws.Range("K61").Value = 787.5
$ws.Range("M61").Value = -615.5
$ws.Range("H62").Value = 5099.7144
$ws.Range("I62").Value = 4551.25
$ws.Range("J62").Value = 5437.231
$ws.Range("K62").Value = 4551.25
$ws.Range("L62").Value = 5437.231
$ws.Range("M62").Value = -3927.25
$ws.Range("N62").Value = -6685.231
$ws.Range("H65").Value = 5099.7144
$ws.Range("I65").Value = 4551.25
$ws.Range("J65").Value = 5437.231
$ws.Range("K65").Value = 22756.25
$ws.Range("L65").Value = 27186.155
$ws.Range("M65").Value = -19636.25
$ws.Range("N65").Value = -33426.155
$ws.Range("H88").Value = 11001.167
$ws.Range("I88").Value = 10001.5
$ws.Range("J88").Value = 11501.0
$ws.Range("K88").Value = 10001.5
$ws.Range("L88").Value = 11501.0
$ws.Range("M88").Value = -9595.5
$ws.Range("N88").Value = -12313.0
$ws.Range("H91").Value = 11001.167
$ws.Range("I91").Value = 10001.5
$ws.Range("J91").Value = 11501.0
$ws.Range("K91").Value = 10001.5
$ws.Range("L91").Value = 11501.0
$ws.Range("M91").Value = -8597.5
$ws.Range("N91").Value = -14309.0
$ws.Range("H136").Value = 38337.332
$ws.Range("J136").Value = 38337.332
$ws.Range("L136").Value = 38337.332
$ws.Range("N136").Value = -48537.332
$ws.Range("H141").Value = 1439.7778
$ws.Range("I141").Value = 1003.8125
$ws.Range("J141").Value = 4927.5
$ws.Range("K141").Value = 3011.4375
$ws.Range("L141").Value = 14782.5
$ws.Range("M141").Value = 2168.5625
$ws.Range("N141").Value = -25142.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3471.08
$ws.Range("I32").Value = 2518.1365
$ws.Range("J32").Value = 10459.333
$ws.Range("K32").Value = 2518.1365
$ws.Range("L32").Value = 10459.333
$ws.Range("M32").Value = -2231.1365
$ws.Range("N32").Value = -11033.333
$ws.Range("H45").Value = 1615.4231
$ws.Range("I45").Value = 1448.4103
$ws.Range("J45").Value = 2116.4614
$ws.Range("K45").Value = 1448.4103
$ws.Range("L45").Value = 2116.4614
$ws.Range("M45").Value = -1071.4103
$ws.Range("N45").Value = -2870.4614
$ws.Range("H61").Value = 1458.1041
$ws.Range("I61").Value = 1329.8438
$ws.Range("J61").Value = 1714.625
$ws.Range("K61").Value = 1329.8438
$ws.Range("L61").Value = 1714.625
$ws.Range("M61").Value = -1117.8438
$ws.Range("N61").Value = -2138.625
$ws.Range("H74").Value = 30277.514
$ws.Range("I74").Value = 41261.4
$ws.Range("J74").Value = 2817.8
$ws.Range("K74").Value = 41261.4
$ws.Range("L74").Value = 2817.8
$ws.Range("M74").Value = -40387.4
$ws.Range("N74").Value = -4565.8
$ws.Range("H77").Value = 30277.514
$ws.Range("I77").Value = 41261.4
$ws.Range("J77").Value = 2817.8
$ws.Range("K77").Value = 206307.0
$ws.Range("L77").Value = 14089.0
$ws.Range("M77").Value = -201939.0
$ws.Range("N77").Value = -22825.0
$ws.Range("H96").Value = 28184.4
$ws.Range("J96").Value = 28184.4
$ws.Range("L96").Value = 28184.4
$ws.Range("N96").Value = -33676.4
$ws.Range("H97").Value = 513.3333
$ws.Range("I97").Value = 533.75
$ws.Range("J97").Value = 350.0
$ws.Range("K97").Value = 533.75
$ws.Range("L97").Value = 350.0
$ws.Range("M97").Value = -37.75
$ws.Range("N97").Value = -1342.0
$ws.Range("H132").Value = 1366.6666
$ws.Range("I132").Value = 1193.8182
$ws.Range("J132").Value = 1638.2858
$ws.Range("K132").Value = 3581.4546
$ws.Range("L132").Value = 4914.857400000001
$ws.Range("M132").Value = -1051.4546
$ws.Range("N132").Value = -9974.8574
$ws.Range("H136").Value = 1458.1041
$ws.Range("I136").Value = 1329.8438
$ws.Range("J136").Value = 1714.625
$ws.Range("K136").Value = 3989.5314
$ws.Range("L136").Value = 5143.875
$ws.Range("M136").Value = -1439.5314
$ws.Range("N136").Value = -10243.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 29115.6
$ws.Range("J55").Value = 29115.6
$ws.Range("L55").Value = 29115.6
$ws.Range("N55").Value = -29661.6
$ws.Range("H111").Value = 24333.0
$ws.Range("J111").Value = 24333.0
$ws.Range("L111").Value = 24333.0
$ws.Range("N111").Value = -32513.0
$ws.Range("H134").Value = 591617.1
$ws.Range("I134").Value = 978434.0
$ws.Range("J134").Value = 4228.5557
$ws.Range("K134").Value = 2935302.0
$ws.Range("L134").Value = 12685.6671
$ws.Range("M134").Value = -2932767.0
$ws.Range("N134").Value = -17755.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 4169.0356
$ws.Range("I25").Value = 2155.5417
$ws.Range("J25").Value = 16250.0
$ws.Range("K25").Value = 2155.5417
$ws.Range("L25").Value = 16250.0
$ws.Range("M25").Value = -1981.5417
$ws.Range("N25").Value = -16598.0
$ws.Range("H31").Value = 8066366.0
$ws.Range("I31").Value = 1148.2609
$ws.Range("J31").Value = 12822776.0
$ws.Range("K31").Value = 1148.2609
$ws.Range("L31").Value = 12822776.0
$ws.Range("M31").Value = -853.2609
$ws.Range("N31").Value = -12823366.0
$ws.Range("H34").Value = 8066366.0
$ws.Range("I34").Value = 1148.2609
$ws.Range("J34").Value = 12822776.0
$ws.Range("K34").Value = 1148.2609
$ws.Range("L34").Value = 12822776.0
$ws.Range("M34").Value = -946.2609
$ws.Range("N34").Value = -12823180.0
$ws.Range("H58").Value = 4054.7576
$ws.Range("I58").Value = 4952.12
$ws.Range("J58").Value = 1250.5
$ws.Range("K58").Value = 4952.12
$ws.Range("L58").Value = 1250.5
$ws.Range("M58").Value = -4749.12
$ws.Range("N58").Value = -1656.5
$ws.Range("H132").Value = 629262.9
$ws.Range("I132").Value = 1458.7551
$ws.Range("J132").Value = 3705503.0
$ws.Range("K132").Value = 4376.2653
$ws.Range("L132").Value = 11116509.0
$ws.Range("M132").Value = -1846.2653
$ws.Range("N132").Value = -11121569.0
$ws.Range("H136").Value = 4054.7576
$ws.Range("I136").Value = 4952.12
$ws.Range("J136").Value = 1250.5
$ws.Range("K136").Value = 14856.36
$ws.Range("L136").Value = 3751.5
$ws.Range("M136").Value = -12306.36
$ws.Range("N136").Value = -8851.5
$ws.Range("H141").Value = 81666.664
$ws.Range("J141").Value = 81666.664
$ws.Range("L141").Value = 81666.664
$ws.Range("N141").Value = -92026.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 25000104.0
$ws.Range("I7").Value = 56.25
$ws.Range("K7").Value = 168.75
$ws.Range("M7").Value = -56.75
$ws.Range("H21").Value = 3323.5881
$ws.Range("I21").Value = 2850.1
$ws.Range("K21").Value = 8550.3
$ws.Range("M21").Value = -8377.3
$ws.Range("H39").Value = 2437.5
$ws.Range("I39").Value = 562.5
$ws.Range("J39").Value = 3375.0
$ws.Range("K39").Value = 1687.5
$ws.Range("L39").Value = 10125.0
$ws.Range("M39").Value = -1393.5
$ws.Range("N39").Value = -10713.0
$ws.Range("H51").Value = 3711.4814
$ws.Range("J51").Value = 3711.4814
$ws.Range("L51").Value = 11134.4442
$ws.Range("N51").Value = -12054.4442
$ws.Range("H109").Value = 1911.9048
$ws.Range("I109").Value = 996.38464
$ws.Range("J109").Value = 3399.625
$ws.Range("K109").Value = 2989.15392
$ws.Range("L109").Value = 10198.875
$ws.Range("M109").Value = -1949.15392
$ws.Range("N109").Value = -12278.875
$ws.Range("H131").Value = 918.7
$ws.Range("J131").Value = 920.81635
$ws.Range("L131").Value = 2762.44905
$ws.Range("N131").Value = -12842.44905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 16427.092
$ws.Range("J57").Value = 17919.8
$ws.Range("L57").Value = 17919.8
$ws.Range("N57").Value = -19559.8
$ws.Range("H122").Value = 70147.266
$ws.Range("I122").Value = 93790.37
$ws.Range("J122").Value = 5128.75
$ws.Range("K122").Value = 281371.11
$ws.Range("L122").Value = 15386.25
$ws.Range("M122").Value = -278921.11
$ws.Range("N122").Value = -20286.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1864.7894
$ws.Range("I136").Value = 1284.8695
$ws.Range("J136").Value = 2754.0
$ws.Range("K136").Value = 3854.6085
$ws.Range("L136").Value = 8262.0
$ws.Range("M136").Value = -1304.6085
$ws.Range("N136").Value = -13362.0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2009.7742
$ws.Range("I132").Value = 2259.0732
$ws.Range("J132").Value = 1523.0476
$ws.Range("K132").Value = 6777.219599999999
$ws.Range("L132").Value = 4569.142800000001
$ws.Range("M132").Value = -4247.219599999999
$ws.Range("N132").Value = -9629.142800000001
